$wb = $excel.ActiveWorkbook

# --- Sheet "Canada": add row 15 -------------------------------------------
$wsCanada = $wb.Worksheets.Item("Canada")

$wsCanada.Cells.Item(15, 1).Value = 44228
$wsCanada.Cells.Item(15, 1).NumberFormat = "d-mmm-yy"
$wsCanada.Cells.Item(15, 2).Value = "Canada"
$wsCanada.Cells.Item(15, 2).NumberFormat = "d-mmm-yy"
$wsCanada.Cells.Item(15, 3).Value = 45.3
$wsCanada.Cells.Item(15, 4).Value = 1665.1

# --- Sheet "Province": add rows 132-141 ------------------------------------
$wsProvince = $wb.Worksheets.Item("Province")

$provinceRows = @(
    @{ Row = 132; Name = "Newfoundland & Labrador"; C = 18;   D = 37.3 },
    @{ Row = 133; Name = "Prince Edward Island";     C = 11.3; D = 7.9 },
    @{ Row = 134; Name = "Nova Scotia";              C = 1;    D = 41 },
    @{ Row = 135; Name = "New Brunswick";            C = 22.4; D = 35 },
    @{ Row = 136; Name = "Quebec";                   C = 41.8; D = 289.4 },
    @{ Row = 137; Name = "Ontario";                  C = 67.1; D = 726.5 },
    @{ Row = 138; Name = "Manitoba";                 C = 30.9; D = 47 },
    @{ Row = 139; Name = "Saskatchewan";             C = 11.9; D = 43.1 },
    @{ Row = 140; Name = "Alberta";                  C = 32.5; D = 242.1 },
    @{ Row = 141; Name = "British Columbia";         C = 38;   D = 195.8 }
)

foreach ($entry in $provinceRows) {
    $r = $entry.Row
    $wsProvince.Cells.Item($r, 1).Value = 44228
    $wsProvince.Cells.Item($r, 1).NumberFormat = "d-mmm-yy"
    $wsProvince.Cells.Item($r, 2).Value = $entry.Name
    $wsProvince.Cells.Item($r, 3).Value = $entry.C
    $wsProvince.Cells.Item($r, 4).Value = $entry.D
}
# First new province row (132) mirrors row 122's formatting: A & B both carry
# the date number format.
$wsProvince.Cells.Item(132, 2).NumberFormat = "d-mmm-yy"

# --- Selections: update the saved cursor position on each sheet -----------
# Visit "Canada" first (sets its stored selection), then finish on
# "Province" so it remains the active/visible tab, matching the workbook's
# original state.
$wsCanada.Range("C16").Select()
$wsProvince.Range("D1").Select()
